{"js": "// Termo de devolu\u00e7\u00e3o: troca o item de lista \"(X) Carregador completo\" pelo\n// placeholder \"{{perifericos}}\" (sem numera\u00e7\u00e3o, com recuo equivalente) e\n// junta os dois runs de \"O mesmo ... fins profissionais;\" (removendo as\n// marcas de proofErr) em um \u00fanico run.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet peripheralsParagraph = null;\nlet sameEquipmentParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text;\n  if (peripheralsParagraph === null && text.indexOf(\"Carregador completo\") !== -1) {\n    peripheralsParagraph = p;\n  }\n  if (sameEquipmentParagraph === null && text.indexOf(\"foi utilizado exclusivamente para fins profissionais\") !== -1) {\n    sameEquipmentParagraph = p;\n  }\n}\n\nif (!peripheralsParagraph) {\n  throw new Error(\"Paragraph with 'Carregador completo' not found\");\n}\nif (!sameEquipmentParagraph) {\n  throw new Error(\"Paragraph with 'foi utilizado exclusivamente para fins profissionais' not found\");\n}\n\n// 1) \"(X) Carregador completo\" -> \"{{perifericos}}\", sem marcador de lista,\n//    com recuo w:ind left=720 no lugar.\nconst peripheralsOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr>\n    <w:ind w:left=\"720\"/>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n    <w:t>{{perifericos}}</w:t>\n  </w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nperipheralsParagraph.getRange().insertOoxml(peripheralsOoxml, Word.InsertLocation.replace);\n\n// 2) \"O mesmo\" + \" foi utilizado exclusivamente para fins profissionais;\"\n//    (dois runs separados por proofErr) -> um \u00fanico run com o texto completo.\nconst sameEquipmentOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr>\n    <w:numPr>\n      <w:ilvl w:val=\"0\"/>\n      <w:numId w:val=\"3\"/>\n    </w:numPr>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n    <w:t>O mesmo foi utilizado exclusivamente para fins profissionais;</w:t>\n  </w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nsameEquipmentParagraph.getRange().insertOoxml(sameEquipmentOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Termo de devolu\u00e7\u00e3o: troca o item de lista \"(X) Carregador completo\" pelo\n# placeholder \"{{perifericos}}\" (sem numera\u00e7\u00e3o, com recuo equivalente) e\n# junta os dois runs de \"O mesmo ... fins profissionais;\" (removendo as\n# marcas de proofErr) em um \u00fanico run.\n\n$d = $word.ActiveDocument\n\n$peripheralsParagraph = $null\n$sameEquipmentParagraph = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if (($null -eq $peripheralsParagraph) -and ($t -like \"*Carregador completo*\")) {\n        $peripheralsParagraph = $p\n    }\n    if (($null -eq $sameEquipmentParagraph) -and ($t -like \"*foi utilizado exclusivamente para fins profissionais*\")) {\n        $sameEquipmentParagraph = $p\n    }\n}\n\nif ($null -eq $peripheralsParagraph) {\n    throw \"Paragraph with 'Carregador completo' not found\"\n}\nif ($null -eq $sameEquipmentParagraph) {\n    throw \"Paragraph with 'foi utilizado exclusivamente para fins profissionais' not found\"\n}\n\n# 1) \"(X) Carregador completo\" -> \"{{perifericos}}\", sem marcador de lista,\n#    com recuo w:ind left=720 no lugar.\n$peripheralsOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr>\n    <w:ind w:left=\"720\"/>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n    <w:t>{{perifericos}}</w:t>\n  </w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$peripheralsParagraph.Range.InsertXML($peripheralsOoxml)\n\n# 2) \"O mesmo\" + \" foi utilizado exclusivamente para fins profissionais;\"\n#    (dois runs separados por proofErr) -> um \u00fanico run com o texto completo.\n$sameEquipmentOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n  <w:pPr>\n    <w:numPr>\n      <w:ilvl w:val=\"0\"/>\n      <w:numId w:val=\"3\"/>\n    </w:numPr>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n    <w:t>O mesmo foi utilizado exclusivamente para fins profissionais;</w:t>\n  </w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$sameEquipmentParagraph.Range.InsertXML($sameEquipmentOoxml)\n"}
